$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Select D2:D12 and type 4, filling the whole selection with the same value
# (D2 was already 4; D3:D12 previously held assorted placeholder numbers)
$ws.Range("D2:D12").Select()
$ws.Range("D2:D12").Value = 4
